$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Adding starting control plot": the sample at index 11 (row 16) becomes
# the new starting control-chart point, column B is widened so the wider
# header/values are readable, and the selection follows the edited cell.

# 1) Update the control-plot starting value in B16 (was 2.18235828).
$ws.Range("B16").Value = -2.51

# 2) Widen column B (was ~11.43 chars) to fit the relabeled control data.
$ws.Columns("B:B").ColumnWidth = 19.1666667

# 3) Nudge the sheet's standard/default column width slightly wider too.
#    (No-op on engines that don't persist StandardWidth, but this is the
#    correct COM call for it.)
$ws.StandardWidth = 11.58984375

# 4) Leave the selection on the cell that was just edited.
$ws.Range("B16").Select()
